# Update the Data Dictionary table (Apartment/Phone rows) to reflect the
# corrected field names, data types, sizes and notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: ApartmentNumber -> Apartment, size 3 -> 4, note "1 per " -> "1 per Apt"
$ws.Range("C12").Value = "Apartment"
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = "1 per Apt"

# Row 13: PhoneNumber -> Phone, type int -> Varchar, size 10 -> 12, add note
$ws.Range("C13").Value = "Phone"
$ws.Range("E13").Value = "Varchar"
$ws.Range("F13").Value = 12
$ws.Range("G13").Value = "(###)###-####"

# Row 14: Email size 30 -> 45
$ws.Range("F14").Value = 45

# Row 15: Password size 20 -> 45
$ws.Range("F15").Value = 45

# Row 17: ScheduleTIme size 10 -> 15
$ws.Range("F17").Value = 15
